# #CRM-31 Remove ID, Bank Details, bracket flag from Download SF list
#
# The "Download SF list" template had columns for the vendor's internal
# ID, bank details (Bank Name / Bank Account / IFSC Code / Beneficiary
# Name) and an internal "Brackets Flag" debug column. These are removed
# from the exported sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds "ID" / "{vendor:id}" - drop it entirely, shifting
# everything after it one column to the left.
$ws.Columns("B").Delete()

# After the shift above, the bank-details block (Bank Name, Bank
# Account, IFSC Code, Beneficiary Name) now lives in columns Z:AC -
# remove all four together.
$ws.Range("Z1:AC1").EntireColumn.Delete()

# Finally, drop the trailing "Brackets Flag" column (now column AB).
$ws.Columns("AB").Delete()

# Restore the sheet's active selection (the author had scrolled right,
# looking at the tail of the shortened row) so the saved view matches.
$ws.Range("AA13").Select() | Out-Null
